$wb = $excel.ActiveWorkbook

# Update F2/F3/F5 (想去人数 / "wish-to-go" attendee counts) on both the
# "展览" and "全部类型" sheets, which carry the same underlying data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1825
    $ws.Range("F3").Value = 8271
    $ws.Range("F5").Value = 333
}
